$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected / new literal numeric values ---
$ws.Range("C2").Value = 0.6163861638616386
$ws.Range("D2").Value = -0.02270816735806249
$ws.Range("E2").Value = 2.541496219221262
$ws.Range("F2").Value = 0.04992291169153087
$ws.Range("C3").Value = 0.6541565415654158
$ws.Range("D3").Value = -0.02647881142360493
$ws.Range("E3").Value = 3.383110910092166
$ws.Range("F3").Value = 0.04641072438264935
$ws.Range("C4").Value = 303.5982859828598
$ws.Range("D4").Value = -0.0000570533119103516
$ws.Range("E4").Value = 39.44031758421771
$ws.Range("F4").Value = 0.003981018349280512
$ws.Range("C5").Value = 20.83175831758317
$ws.Range("D5").Value = -0.006352483016761816
$ws.Range("C6").Value = 10.62520625206252
$ws.Range("D6").Value = -0.01245466561141303
$ws.Range("C7").Value = 2.8750037500375
$ws.Range("D7").Value = -0.0463148395710653
$ws.Range("E7").Value = 3.378140312826196
$ws.Range("F7").Value = 0.357301881058571
$ws.Range("C8").Value = 2.985554855548555
$ws.Range("D8").Value = -0.04459986297077515
$ws.Range("E8").Value = 3.379908351634691
$ws.Range("F8").Value = 0.3571149755196222
$ws.Range("C9").Value = 0.3775462754627545
$ws.Range("D9").Value = -0.3172204919802412
$ws.Range("B12").Value = 198.0
$ws.Range("D12").Value = 2.0569
$ws.Range("B13").Value = 160.0
$ws.Range("D13").Value = 2.0569
$ws.Range("B14").Value = 160.0
$ws.Range("D14").Value = 2.0569
$ws.Range("B15").Value = 162.0
$ws.Range("D15").Value = 15.911
$ws.Range("B16").Value = 161.0
$ws.Range("D16").Value = 15.911
$ws.Range("B17").Value = 179.0
$ws.Range("D17").Value = 15.911

# --- New literal string values ---
$ws.Range("A12").Value = "short period 1"
$ws.Range("A13").Value = "short period 2"
$ws.Range("A14").Value = "phugoid"
$ws.Range("A15").Value = "spiral"
$ws.Range("A16").Value = "dutch roll"
$ws.Range("A17").Value = "roll damping"
$ws.Range("B11").Value = "V kts"
$ws.Range("C11").Value = "V m/s"
$ws.Range("D11").Value = "mac/b"
$ws.Range("E11").Value = "(mac/b)/V"

# --- Formulas (values computed automatically by the engine) ---
$ws.Range("H2").Formula = "=D2/E12"
$ws.Range("J2").Formula = "=F2/E12"
$ws.Range("H3").Formula = "=D3/E13"
$ws.Range("J3").Formula = "=F3/E13"
$ws.Range("H4").Formula = "=D4/E14"
$ws.Range("J4").Formula = "=F4/E14"
$ws.Range("H5").Formula = "=D5/E15"
$ws.Range("H6").Formula = "=D6/E16"
$ws.Range("H7").Formula = "=D7/E16"
$ws.Range("J7").Formula = "=F7/E16"
$ws.Range("H8").Formula = "=D8/E16"
$ws.Range("J8").Formula = "=F8/E16"
$ws.Range("H9").Formula = "=D9/E17"
$ws.Range("C12").Formula = "=B12*0.514444444"
$ws.Range("E12").Formula = "=D12/C12"
$ws.Range("C13").Formula = "=B13*0.514444444"
$ws.Range("E13").Formula = "=D13/C13"
$ws.Range("C14").Formula = "=B14*0.514444444"
$ws.Range("E14").Formula = "=D14/C14"
$ws.Range("C15").Formula = "=B15*0.514444444"
$ws.Range("E15").Formula = "=D15/C15"
$ws.Range("C16").Formula = "=B16*0.514444444"
$ws.Range("E16").Formula = "=D16/C16"
$ws.Range("C17").Formula = "=B17*0.514444444"
$ws.Range("E17").Formula = "=D17/C17"

# --- Column widths ---
# (ColumnWidth is rounded to the nearest 1/6 character-unit by the
#  engine before storing OOXML width = ColumnWidth + 5/6, so the values
#  below are chosen to land as close as possible to the target widths
#  13.42578125 / 17.28515625 / 21.7109375 / 21.85546875.)
$ws.Columns.Item(1).ColumnWidth = 12.666666666666666
$ws.Columns.Item(2).ColumnWidth = 16.5
$ws.Columns.Item(4).ColumnWidth = 20.833333333333332
$ws.Columns.Item(6).ColumnWidth = 21.0

# --- Selection / view state ---
$ws.Range("H3").Select()
